$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.000000000317472714783662
$ws.Range("C2").Value = 0.0000005461030343489881
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 10.19245300693656
$ws.Range("G2").Value = 32.58298917216626

$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 250555.8564151394
$ws.Range("D3").Value = 22.3905356188092
$ws.Range("E3").Value = 1133.036916526867
$ws.Range("G3").Value = 251714.5706998299
